$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before current column B -> becomes "Region"
#    (Date/CFR2/DiffSK/AgeCompSK/RateCompSK/relAgeSK/relRateSK shift from B..H to C..I)
$ws.Columns.Item(2).Insert()

# 2) Insert a new row before current row 2 -> new first data row (USA / All)
$ws.Rows.Item(2).Insert()

# Reset the formatting of the freshly inserted row so it does not inherit
# the bold/centered header style; bring it back to the workbook default.
$ws.Range("A2:I2").Style = "Normal"

# Re-apply the date number format on the new row's Date cell (column C)
$ws.Range("C2").NumberFormat = "yyyy-mm-dd"

# 3) Header row
$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Region"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "CFR2"
$ws.Range("E1").Value = "DiffSK"
$ws.Range("F1").Value = "AgeCompSK"
$ws.Range("G1").Value = "RateCompSK"
$ws.Range("H1").Value = "relAgeSK"
$ws.Range("I1").Value = "relRateSK"

# 4) Data rows - (Country, Region, Date, CFR2, DiffSK, AgeCompSK, RateCompSK, relAgeSK, relRateSK)
$rows = @(
  @("USA",        "All", 43906, 0.01796725060231124,  0.003981909889199872,   -0.01231536097948927,   0.01629727086868914,   0.4304169237152266,  0.5695830762847735),
  @("SouthKorea",  "All", 43943, 0.02225547035720965,  -0.0003063098656985348, 0.00001981077847107633, -0.0003261206441696108, 0.05726793570774701, 0.942732064292253),
  @("China",       "All", 43872, 0.02290025071633238,  -0.000951090224821264,  -0.001865905070842499,  0.0009148148460212372,  0.6710151063854642,  0.3289848936145358),
  @("Germany",     "All", 43941, 0.03205143844597228,  -0.01010227795446117,   -0.01224049443597415,   0.002138216481512984,   0.8512928944894133,  0.1487071055105867),
  @("France",      "All", 43914, 0.03983587515221891,  -0.0178867146607078,    -0.02067143833938045,   0.002784723678672646,   0.8812796536564944,  0.1187203463435056),
  @("USA",         "NYC", 43941, 0.0708957990420689,   -0.04894663855055778,   -0.01281025699693983,   -0.03613638155361795,   0.2617188304710222,  0.7382811695289777),
  @("Spain",       "All", 43937, 0.1050210003716739,   -0.08307183988016276,   -0.05560002378836659,   -0.02747181609179617,   0.66930049784107,    0.33069950215893),
  @("Italy",       "All", 43941, 0.1272752828730058,   -0.1053261223814947,    -0.06796855464528848,   -0.03735756773620624,   0.6453152656574986,  0.3546847343425013)
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $ws.Cells.Item($r, 9).Value = $row[8]
  $r = $r + 1
}
